$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("G2").Value = 55.59510733333332
$ws.Range("H2").Value = 166.785322
$ws.Range("I2").Value = 0.4537221086682116
$ws.Range("J2").Value = 0.4537221086682116
$ws.Range("K2").Value = 3
$ws.Range("M2").Value = 1.655851666666667
$ws.Range("N2").Value = 4.967555
$ws.Range("O2").Value = 0.03628213169899143
$ws.Range("P2").Value = 0.03628213169899143
$ws.Range("Q2").Value = 92.05725113641221
$ws.Range("R2").Value = 828.5152602277099
$ws.Range("S2").Value = 0.01646200530144416
$ws.Range("T2").Value = 0.01646200530144415
$ws.Range("E3").Value = 3
$ws.Range("G3").Value = 55.59510733333332
$ws.Range("H3").Value = 166.785322
$ws.Range("I3").Value = 0.4537221086682116
$ws.Range("J3").Value = 0.4537221086682116
$ws.Range("K3").Value = 3
$ws.Range("M3").Value = 1.416382666666667
$ws.Range("N3").Value = 4.249148
$ws.Range("O3").Value = 0.03103501568568562
$ws.Range("P3").Value = 0.03103501568568562
$ws.Range("Q3").Value = 78.74394637840621
$ws.Range("R3").Value = 708.6955174056559
$ws.Range("S3").Value = 0.0140812727594603
$ws.Range("T3").Value = 0.0140812727594603
$ws.Range("E4").Value = 3
$ws.Range("G4").Value = 55.59510733333332
$ws.Range("H4").Value = 166.785322
$ws.Range("I4").Value = 0.4537221086682116
$ws.Range("J4").Value = 0.4537221086682116
$ws.Range("K4").Value = 3
$ws.Range("M4").Value = 5.385314999999999
$ws.Range("N4").Value = 16.155945
$ws.Range("O4").Value = 0.1180001276707882
$ws.Range("P4").Value = 0.1180001276707882
$ws.Range("Q4").Value = 299.3971654488099
$ws.Range("R4").Value = 2694.574489039289
$ws.Range("S4").Value = 0.05353926674990819
$ws.Range("T4").Value = 0.05353926674990819
$ws.Range("E5").Value = 3
$ws.Range("G5").Value = 55.59510733333332
$ws.Range("H5").Value = 166.785322
$ws.Range("I5").Value = 0.4537221086682116
$ws.Range("J5").Value = 0.4537221086682116
$ws.Range("K5").Value = 3
$ws.Range("M5").Value = 37.180664
$ws.Range("N5").Value = 111.541992
$ws.Range("O5").Value = 0.8146827249445348
$ws.Range("P5").Value = 0.8146827249445348
$ws.Range("Q5").Value = 2067.063005804602
$ws.Range("R5").Value = 18603.56705224142
$ws.Range("S5").Value = 0.3696395638573989
$ws.Range("T5").Value = 0.3696395638573989
$ws.Range("E6").Value = 3
$ws.Range("G6").Value = 1.654659
$ws.Range("H6").Value = 4.963977
$ws.Range("I6").Value = 0.01350398275347337
$ws.Range("J6").Value = 0.01350398275347337
$ws.Range("K6").Value = 3
$ws.Range("M6").Value = 1.655851666666667
$ws.Range("N6").Value = 4.967555
$ws.Range("O6").Value = 0.03628213169899143
$ws.Range("P6").Value = 0.03628213169899143
$ws.Range("Q6").Value = 2.739869862915
$ws.Range("R6").Value = 24.658828766235
$ws.Range("S6").Value = 0.0004899532807224299
$ws.Range("T6").Value = 0.0004899532807224298
$ws.Range("E7").Value = 3
$ws.Range("G7").Value = 1.654659
$ws.Range("H7").Value = 4.963977
$ws.Range("I7").Value = 0.01350398275347337
$ws.Range("J7").Value = 0.01350398275347337
$ws.Range("K7").Value = 3
$ws.Range("M7").Value = 1.416382666666667
$ws.Range("N7").Value = 4.249148
$ws.Range("O7").Value = 0.03103501568568562
$ws.Range("P7").Value = 0.03103501568568562
$ws.Range("Q7").Value = 2.343630326844
$ws.Range("R7").Value = 21.092672941596
$ws.Range("S7").Value = 0.0004190963165732742
$ws.Range("T7").Value = 0.0004190963165732742
$ws.Range("E8").Value = 3
$ws.Range("G8").Value = 1.654659
$ws.Range("H8").Value = 4.963977
$ws.Range("I8").Value = 0.01350398275347337
$ws.Range("J8").Value = 0.01350398275347337
$ws.Range("K8").Value = 3
$ws.Range("M8").Value = 5.385314999999999
$ws.Range("N8").Value = 16.155945
$ws.Range("O8").Value = 0.1180001276707882
$ws.Range("P8").Value = 0.1180001276707882
$ws.Range("Q8").Value = 8.910859932584998
$ws.Range("R8").Value = 80.197739393265
$ws.Range("S8").Value = 0.001593471688973979
$ws.Range("T8").Value = 0.001593471688973979
$ws.Range("E9").Value = 3
$ws.Range("G9").Value = 1.654659
$ws.Range("H9").Value = 4.963977
$ws.Range("I9").Value = 0.01350398275347337
$ws.Range("J9").Value = 0.01350398275347337
$ws.Range("K9").Value = 3
$ws.Range("M9").Value = 37.180664
$ws.Range("N9").Value = 111.541992
$ws.Range("O9").Value = 0.8146827249445348
$ws.Range("P9").Value = 0.8146827249445348
$ws.Range("Q9").Value = 61.521320313576
$ws.Range("R9").Value = 553.691882822184
$ws.Range("S9").Value = 0.01100146146720369
$ws.Range("T9").Value = 0.01100146146720369
$ws.Range("E10").Value = 3
$ws.Range("G10").Value = 62.05924166666666
$ws.Range("H10").Value = 186.177725
$ws.Range("I10").Value = 0.5064771225734745
$ws.Range("J10").Value = 0.5064771225734744
$ws.Range("K10").Value = 3
$ws.Range("M10").Value = 1.655851666666667
$ws.Range("N10").Value = 4.967555
$ws.Range("O10").Value = 0.03628213169899143
$ws.Range("P10").Value = 0.03628213169899143
$ws.Range("Q10").Value = 102.7608987458194
$ws.Range("R10").Value = 924.8480887123749
$ws.Range("S10").Value = 0.01837606966373703
$ws.Range("T10").Value = 0.01837606966373703
$ws.Range("E11").Value = 3
$ws.Range("G11").Value = 62.05924166666666
$ws.Range("H11").Value = 186.177725
$ws.Range("I11").Value = 0.5064771225734745
$ws.Range("J11").Value = 0.5064771225734744
$ws.Range("K11").Value = 3
$ws.Range("M11").Value = 1.416382666666667
$ws.Range("N11").Value = 4.249148
$ws.Range("O11").Value = 0.03103501568568562
$ws.Range("P11").Value = 0.03103501568568562
$ws.Range("Q11").Value = 87.89963420314443
$ws.Range("R11").Value = 791.0967078282999
$ws.Range("S11").Value = 0.0157185254435087
$ws.Range("T11").Value = 0.0157185254435087
$ws.Range("E12").Value = 3
$ws.Range("G12").Value = 62.05924166666666
$ws.Range("H12").Value = 186.177725
$ws.Range("I12").Value = 0.5064771225734745
$ws.Range("J12").Value = 0.5064771225734744
$ws.Range("K12").Value = 3
$ws.Range("M12").Value = 5.385314999999999
$ws.Range("N12").Value = 16.155945
$ws.Range("O12").Value = 0.1180001276707882
$ws.Range("P12").Value = 0.1180001276707882
$ws.Range("Q12").Value = 334.2085650361249
$ws.Range("R12").Value = 3007.877085325124
$ws.Range("S12").Value = 0.05976436512600342
$ws.Range("T12").Value = 0.05976436512600342
$ws.Range("E13").Value = 3
$ws.Range("G13").Value = 62.05924166666666
$ws.Range("H13").Value = 186.177725
$ws.Range("I13").Value = 0.5064771225734745
$ws.Range("J13").Value = 0.5064771225734744
$ws.Range("K13").Value = 3
$ws.Range("M13").Value = 37.180664
$ws.Range("N13").Value = 111.541992
$ws.Range("O13").Value = 0.8146827249445348
$ws.Range("P13").Value = 0.8146827249445348
$ws.Range("Q13").Value = 2307.403812503133
$ws.Range("R13").Value = 20766.6343125282
$ws.Range("S13").Value = 0.4126181623402254
$ws.Range("T13").Value = 0.4126181623402253
$ws.Range("E14").Value = 3
$ws.Range("G14").Value = 3.222176333333334
$ws.Range("H14").Value = 9.666529000000001
$ws.Range("I14").Value = 0.02629678600484052
$ws.Range("J14").Value = 0.02629678600484052
$ws.Range("K14").Value = 3
$ws.Range("M14").Value = 1.655851666666667
$ws.Range("N14").Value = 4.967555
$ws.Range("O14").Value = 0.03628213169899143
$ws.Range("P14").Value = 0.03628213169899143
$ws.Range("Q14").Value = 5.335446051843889
$ws.Range("R14").Value = 48.019014466595
$ws.Range("S14").Value = 0.0009541034530878184
$ws.Range("T14").Value = 0.0009541034530878182
$ws.Range("E15").Value = 3
$ws.Range("G15").Value = 3.222176333333334
$ws.Range("H15").Value = 9.666529000000001
$ws.Range("I15").Value = 0.02629678600484052
$ws.Range("J15").Value = 0.02629678600484052
$ws.Range("K15").Value = 3
$ws.Range("M15").Value = 1.416382666666667
$ws.Range("N15").Value = 4.249148
$ws.Range("O15").Value = 0.03103501568568562
$ws.Range("P15").Value = 0.03103501568568562
$ws.Range("Q15").Value = 4.563834707476889
$ws.Range("R15").Value = 41.074512367292
$ws.Range("S15").Value = 0.0008161211661433436
$ws.Range("T15").Value = 0.0008161211661433436
$ws.Range("E16").Value = 3
$ws.Range("G16").Value = 3.222176333333334
$ws.Range("H16").Value = 9.666529000000001
$ws.Range("I16").Value = 0.02629678600484052
$ws.Range("J16").Value = 0.02629678600484052
$ws.Range("K16").Value = 3
$ws.Range("M16").Value = 5.385314999999999
$ws.Range("N16").Value = 16.155945
$ws.Range("O16").Value = 0.1180001276707882
$ws.Range("P16").Value = 0.1180001276707882
$ws.Range("Q16").Value = 17.352434540545
$ws.Range("R16").Value = 156.171910864905
$ws.Range("S16").Value = 0.003103024105902576
$ws.Range("T16").Value = 0.003103024105902576
$ws.Range("E17").Value = 3
$ws.Range("G17").Value = 3.222176333333334
$ws.Range("H17").Value = 9.666529000000001
$ws.Range("I17").Value = 0.02629678600484052
$ws.Range("J17").Value = 0.02629678600484052
$ws.Range("K17").Value = 3
$ws.Range("M17").Value = 37.180664
$ws.Range("N17").Value = 111.541992
$ws.Range("O17").Value = 0.8146827249445348
$ws.Range("P17").Value = 0.8146827249445348
$ws.Range("Q17").Value = 119.8026555984187
$ws.Range("R17").Value = 1078.223900385768
$ws.Range("S17").Value = 0.02142353727970678
$ws.Range("T17").Value = 0.02142353727970678
